$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new shared-string content: prepend "Boolean" and "String" entries to N2's existing
# value list, and add a new token list value "[3.1; double]" to the new Q3 cell.

$n2 = $ws.Range("N2")
$n2.Value = "[Boolean; java.lang], [String; java.lang], " + $n2.Value2

$ws.Range("Q3").Value = "[3.1; double]"

# Row 3 height adjustment
$ws.Rows.Item(3).RowHeight = 409.5

# Update the view: scrolled position and active selection
$excel.ActiveWindow.ScrollColumn = 16
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("R3").Select()
